$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: restructure rows 37-39 before filling data ---
# Move the trailing blank row (was row 37) down to its new location row 39
$ws.Range("A37:F37").Copy($ws.Range("A39:F39"))
# Build out two new data rows (37, 38) with the correct alternating style
# (row 37 odd -> style like row 9; row 38 even -> style like row 8)
$ws.Range("A9:F9").Copy($ws.Range("A37:F37"))
$ws.Range("A8:F8").Copy($ws.Range("A38:F38"))

# --- Step 2: write the new cell values as TEXT, preserving each rows existing style ---
# A staging range (used further down the sheet, deleted at the end) lets us force
# every written value to stay a text cell (matching the source data convention),
# then PasteSpecial(xlPasteValues) drops the text onto the destination without
# disturbing the destinations own number format / style / borders.
$stage = $ws.Range("A100:F100")
$stage.NumberFormat = "@"

$ws.Range("A100").Value = "2024-08-31"
$ws.Range("B100").Value = "토"
$ws.Range("C100").Value = "27"
$ws.Range("D100").Value = "1"
$ws.Range("E100").Value = "0"
$ws.Range("F100").Value = "26"
$stage.Copy()
$ws.Range("A8:F8").PasteSpecial(-4163)

$ws.Range("A100").Value = "2024-08-30"
$ws.Range("B100").Value = "금"
$ws.Range("C100").Value = "17"
$ws.Range("D100").Value = "0"
$ws.Range("E100").Value = "0"
$ws.Range("F100").Value = "17"
$stage.Copy()
$ws.Range("A9:F9").PasteSpecial(-4163)

$ws.Range("A100").Value = "2024-08-29"
$ws.Range("B100").Value = "목"
$ws.Range("C100").Value = "23"
$ws.Range("D100").Value = "0"
$ws.Range("E100").Value = "0"
$ws.Range("F100").Value = "23"
$stage.Copy()
$ws.Range("A10:F10").PasteSpecial(-4163)

$ws.Range("A100").Value = "2024-08-28"
$ws.Range("B100").Value = "수"
$ws.Range("C100").Value = "22"
$ws.Range("D100").Value = "0"
$ws.Range("E100").Value = "0"
$ws.Range("F100").Value = "22"
$stage.Copy()
$ws.Range("A11:F11").PasteSpecial(-4163)

$ws.Range("A100").Value = "2024-08-27"
$ws.Range("B100").Value = "화"
$ws.Range("C100").Value = "26"
$ws.Range("D100").Value = "0"
$ws.Range("E100").Value = "0"
$ws.Range("F100").Value = "26"
$stage.Copy()
$ws.Range("A12:F12").PasteSpecial(-4163)

$ws.Range("A100").Value = "2024-08-26"
$ws.Range("B100").Value = "월"
$ws.Range("C100").Value = "25"
$ws.Range("D100").Value = "0"
$ws.Range("E100").Value = "0"
$ws.Range("F100").Value = "25"
$stage.Copy()
$ws.Range("A13:F13").PasteSpecial(-4163)

$ws.Range("A100").Value = "2024-08-25"
$ws.Range("B100").Value = "일"
$ws.Range("C100").Value = "12"
$ws.Range("D100").Value = "0"
$ws.Range("E100").Value = "0"
$ws.Range("F100").Value = "12"
$stage.Copy()
$ws.Range("A14:F14").PasteSpecial(-4163)

$ws.Range("A100").Value = "2024-08-24"
$ws.Range("B100").Value = "토"
$ws.Range("C100").Value = "24"
$ws.Range("D100").Value = "0"
$ws.Range("E100").Value = "0"
$ws.Range("F100").Value = "24"
$stage.Copy()
$ws.Range("A15:F15").PasteSpecial(-4163)

$ws.Range("A100").Value = "2024-08-23"
$ws.Range("B100").Value = "금"
$ws.Range("C100").Value = "20"
$ws.Range("D100").Value = "0"
$ws.Range("E100").Value = "0"
$ws.Range("F100").Value = "20"
$stage.Copy()
$ws.Range("A16:F16").PasteSpecial(-4163)

$ws.Range("A100").Value = "2024-08-22"
$ws.Range("B100").Value = "목"
$ws.Range("C100").Value = "12"
$ws.Range("D100").Value = "0"
$ws.Range("E100").Value = "0"
$ws.Range("F100").Value = "12"
$stage.Copy()
$ws.Range("A17:F17").PasteSpecial(-4163)

$ws.Range("A100").Value = "2024-08-21"
$ws.Range("B100").Value = "수"
$ws.Range("C100").Value = "25"
$ws.Range("D100").Value = "3"
$ws.Range("E100").Value = "0"
$ws.Range("F100").Value = "22"
$stage.Copy()
$ws.Range("A18:F18").PasteSpecial(-4163)

$ws.Range("A100").Value = "2024-08-20"
$ws.Range("B100").Value = "화"
$ws.Range("C100").Value = "21"
$ws.Range("D100").Value = "0"
$ws.Range("E100").Value = "0"
$ws.Range("F100").Value = "21"
$stage.Copy()
$ws.Range("A19:F19").PasteSpecial(-4163)

$ws.Range("A100").Value = "2024-08-19"
$ws.Range("B100").Value = "월"
$ws.Range("C100").Value = "26"
$ws.Range("D100").Value = "2"
$ws.Range("E100").Value = "0"
$ws.Range("F100").Value = "24"
$stage.Copy()
$ws.Range("A20:F20").PasteSpecial(-4163)

$ws.Range("A100").Value = "2024-08-18"
$ws.Range("B100").Value = "일"
$ws.Range("C100").Value = "18"
$ws.Range("D100").Value = "0"
$ws.Range("E100").Value = "0"
$ws.Range("F100").Value = "18"
$stage.Copy()
$ws.Range("A21:F21").PasteSpecial(-4163)

$ws.Range("A100").Value = "2024-08-17"
$ws.Range("B100").Value = "토"
$ws.Range("C100").Value = "11"
$ws.Range("D100").Value = "2"
$ws.Range("E100").Value = "0"
$ws.Range("F100").Value = "9"
$stage.Copy()
$ws.Range("A22:F22").PasteSpecial(-4163)

$ws.Range("A100").Value = "2024-08-16"
$ws.Range("B100").Value = "금"
$ws.Range("C100").Value = "24"
$ws.Range("D100").Value = "0"
$ws.Range("E100").Value = "0"
$ws.Range("F100").Value = "24"
$stage.Copy()
$ws.Range("A23:F23").PasteSpecial(-4163)

$ws.Range("A100").Value = "2024-08-15"
$ws.Range("B100").Value = "목"
$ws.Range("C100").Value = "23"
$ws.Range("D100").Value = "1"
$ws.Range("E100").Value = "0"
$ws.Range("F100").Value = "22"
$stage.Copy()
$ws.Range("A24:F24").PasteSpecial(-4163)

$ws.Range("A100").Value = "2024-08-14"
$ws.Range("B100").Value = "수"
$ws.Range("C100").Value = "24"
$ws.Range("D100").Value = "6"
$ws.Range("E100").Value = "0"
$ws.Range("F100").Value = "18"
$stage.Copy()
$ws.Range("A25:F25").PasteSpecial(-4163)

$ws.Range("A100").Value = "2024-08-13"
$ws.Range("B100").Value = "화"
$ws.Range("C100").Value = "28"
$ws.Range("D100").Value = "0"
$ws.Range("E100").Value = "0"
$ws.Range("F100").Value = "28"
$stage.Copy()
$ws.Range("A26:F26").PasteSpecial(-4163)

$ws.Range("A100").Value = "2024-08-12"
$ws.Range("B100").Value = "월"
$ws.Range("C100").Value = "21"
$ws.Range("D100").Value = "0"
$ws.Range("E100").Value = "0"
$ws.Range("F100").Value = "21"
$stage.Copy()
$ws.Range("A27:F27").PasteSpecial(-4163)

$ws.Range("A100").Value = "2024-08-11"
$ws.Range("B100").Value = "일"
$ws.Range("C100").Value = "18"
$ws.Range("D100").Value = "0"
$ws.Range("E100").Value = "0"
$ws.Range("F100").Value = "18"
$stage.Copy()
$ws.Range("A28:F28").PasteSpecial(-4163)

$ws.Range("A100").Value = "2024-08-10"
$ws.Range("B100").Value = "토"
$ws.Range("C100").Value = "14"
$ws.Range("D100").Value = "0"
$ws.Range("E100").Value = "0"
$ws.Range("F100").Value = "14"
$stage.Copy()
$ws.Range("A29:F29").PasteSpecial(-4163)

$ws.Range("A100").Value = "2024-08-09"
$ws.Range("B100").Value = "금"
$ws.Range("C100").Value = "23"
$ws.Range("D100").Value = "0"
$ws.Range("E100").Value = "0"
$ws.Range("F100").Value = "23"
$stage.Copy()
$ws.Range("A30:F30").PasteSpecial(-4163)

$ws.Range("A100").Value = "2024-08-08"
$ws.Range("B100").Value = "목"
$ws.Range("C100").Value = "20"
$ws.Range("D100").Value = "0"
$ws.Range("E100").Value = "0"
$ws.Range("F100").Value = "20"
$stage.Copy()
$ws.Range("A31:F31").PasteSpecial(-4163)

$ws.Range("A100").Value = "2024-08-07"
$ws.Range("B100").Value = "수"
$ws.Range("C100").Value = "14"
$ws.Range("D100").Value = "9"
$ws.Range("E100").Value = "0"
$ws.Range("F100").Value = "5"
$stage.Copy()
$ws.Range("A32:F32").PasteSpecial(-4163)

$ws.Range("A100").Value = "2024-08-06"
$ws.Range("B100").Value = "화"
$ws.Range("C100").Value = "16"
$ws.Range("D100").Value = "2"
$ws.Range("E100").Value = "0"
$ws.Range("F100").Value = "14"
$stage.Copy()
$ws.Range("A33:F33").PasteSpecial(-4163)

$ws.Range("A100").Value = "2024-08-05"
$ws.Range("B100").Value = "월"
$ws.Range("C100").Value = "14"
$ws.Range("D100").Value = "0"
$ws.Range("E100").Value = "0"
$ws.Range("F100").Value = "14"
$stage.Copy()
$ws.Range("A34:F34").PasteSpecial(-4163)

$ws.Range("A100").Value = "2024-08-04"
$ws.Range("B100").Value = "일"
$ws.Range("C100").Value = "35"
$ws.Range("D100").Value = "1"
$ws.Range("E100").Value = "0"
$ws.Range("F100").Value = "34"
$stage.Copy()
$ws.Range("A35:F35").PasteSpecial(-4163)

$ws.Range("A100").Value = "2024-08-03"
$ws.Range("B100").Value = "토"
$ws.Range("C100").Value = "13"
$ws.Range("D100").Value = "0"
$ws.Range("E100").Value = "0"
$ws.Range("F100").Value = "13"
$stage.Copy()
$ws.Range("A36:F36").PasteSpecial(-4163)

$ws.Range("A100").Value = "2024-08-02"
$ws.Range("B100").Value = "금"
$ws.Range("C100").Value = "20"
$ws.Range("D100").Value = "0"
$ws.Range("E100").Value = "0"
$ws.Range("F100").Value = "20"
$stage.Copy()
$ws.Range("A37:F37").PasteSpecial(-4163)

$ws.Range("A100").Value = "2024-08-01"
$ws.Range("B100").Value = "목"
$ws.Range("C100").Value = "22"
$ws.Range("D100").Value = "0"
$ws.Range("E100").Value = "0"
$ws.Range("F100").Value = "22"
$stage.Copy()
$ws.Range("A38:F38").PasteSpecial(-4163)

# Clean up the staging row entirely so it leaves no trace in the saved file
$ws.Rows("100:100").Delete()

# --- Step 3: update the summary header cells ---
$hstage = $ws.Range("H1")
$hstage.NumberFormat = "@"
$ws.Range("H1").Value = "2024-08-01 ~ 2024-08-31"
$ws.Range("H1").Copy()
$ws.Range("B4").PasteSpecial(-4163)

$ws.Range("H1").Value = "2024년 09월 04일 19시 03분 29초"
$ws.Range("H1").Copy()
$ws.Range("B5").PasteSpecial(-4163)

$ws.Range("H1").ClearContents()
$ws.Range("H1").NumberFormat = "General"